$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "on retire lalheine a soy par le petit pertuis"
#   -> "on retire l" + "<bp>" + "alheine" + "</bp>" + " a soy par le petit pertuis"
# The new "<bp>" / "</bp>" runs must carry the same formatting already used
# elsewhere in the document for inline tag markup (Courier New, blue, 18 half-
# points). We grab that formatting by copying an existing tag run ("<m>") and
# pasting it in place, then swapping its text for the tag we actually want -
# this keeps every w:rPr attribute (including w:eastAsia / w:cs / szCs and the
# lower-case hex color) byte-identical to the rest of the document.
# ---------------------------------------------------------------------------
$tmpl = $d.Content
$tmpl.Find.Execute("<m>") | Out-Null
$tmpl.Copy()

$target = $d.Content
$target.Find.Execute("alheine") | Out-Null
$alheineStart = $target.Start
$alheineEnd = $target.End

# Insert the closing tag right after "alheine" first so that $alheineStart
# (the position right before "alheine") is unaffected by the insertion.
$afterPoint = $d.Range($alheineEnd, $alheineEnd)
$afterPoint.Paste()
$closeTagRange = $d.Range($alheineEnd, $alheineEnd + 3)
$closeTagRange.Text = "</bp>"

# Insert the opening tag right before "alheine".
$beforePoint = $d.Range($alheineStart, $alheineStart)
$beforePoint.Paste()
$openTagRange = $d.Range($alheineStart, $alheineStart + 3)
$openTagRange.Text = "<bp>"

# ---------------------------------------------------------------------------
# Change 2: left-top -> left-middle
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("left-top", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "left-middle", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "Noublye pas " + "dhuiler" -> "Noublye pas d" + "huiler"
#   (the "d" moves from the start of the tagged word to the end of the
#   preceding plain-text run; formatting of each run is untouched otherwise)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Noublye pas ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Noublye pas d", 2) | Out-Null
$d.Content.Find.Execute("dhuiler", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "huiler", 2) | Out-Null
